# Update automatico via Actualizar 02-18-2021 12-42-07
#
# The "Fecha" column (D) holds a rolling history of timestamps in blocks of
# 14 rows each (one block per monitored service group). Each time the
# checker runs, the newest timestamp block shifts down to become the
# "previous" block, and so on; the oldest block is dropped off the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newest = 44245.52911925653
$middle = 44245.50784878472
$oldest = 44245.48657256945

$ws.Range("D2:D15").Value = $newest
$ws.Range("D16:D29").Value = $middle
$ws.Range("D30:D43").Value = $oldest
